$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Internal view Widget" block (old rows 14-21) down to rows 22-29,
# opening up rows 12-13 for new content.
$ws.Range("A14:C21").Insert()

# New row 12: a plain data row.
$ws.Range("A12").Value = "Country"
$ws.Range("B12").Value = "Download button"

# New row 13: data row whose C cell becomes a hyperlink.
$ws.Range("A13").Value = "Engagements"
$ws.Range("B13").Value = "percentages of youths completing each module and dissaggregate by company"

$ws.Hyperlinks.Add(
    $ws.Range("C13"),
    "https://data.yes4youth.co.za/GlobalOverallLearning?csv=1",
    [Type]::Missing,
    [Type]::Missing,
    "https://data.yes4youth.co.za/GlobalOverallLearning?csv=1"
) | Out-Null

# The cell text itself carries the extra "{george email}" annotation, while
# the hyperlink's display/screen text stays the bare URL (set above).
$ws.Range("C13").Value = "(https://data.yes4youth.co.za/GlobalOverallLearning?csv=1) {george email}"

# Match the author's final view/selection state.
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
